$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet: "F-SW-FR-03" -> "F-SW-SD-03" ---
$ws.Name = "F-SW-SD-03"

# Re-apply the print area on the (now renamed) sheet so the workbook-level
# "_xlnm.Print_Area" defined name is regenerated against the new sheet name.
$ws.PageSetup.PrintArea = "A1:D23"

# --- View: zoom to 80% (persists to sheetView zoomScale) ---
$excel.ActiveWindow.Zoom = 80

# --- Footer text update ---
# Old: Rev: 0(0/0/2025)   New: Rev:0(01/10/2025)
$ws.PageSetup.LeftFooter   = '&"Arial,Regular"&14Issue No.:(01)'
$ws.PageSetup.CenterFooter = '&"Arial,Regular"&14F-SW-SD/03'
$ws.PageSetup.RightFooter  = '&"Arial,Regular"&14Rev:0(01/10/2025)'
